$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "91÷7="
$t.Cell(1, 2).Range.Text = "68÷2="
$t.Cell(1, 3).Range.Text = "59÷2="
$t.Cell(1, 4).Range.Text = "17÷9="
$t.Cell(1, 5).Range.Text = "24÷8="

$t.Cell(5, 1).Range.Text = "79÷4="
$t.Cell(5, 2).Range.Text = "80÷8="
$t.Cell(5, 3).Range.Text = "16÷4="
$t.Cell(5, 4).Range.Text = "37÷9="
$t.Cell(5, 5).Range.Text = "34÷8="

$t.Cell(9, 1).Range.Text = "61÷3="
$t.Cell(9, 2).Range.Text = "61÷7="
$t.Cell(9, 3).Range.Text = "93÷2="
$t.Cell(9, 4).Range.Text = "41÷2="
$t.Cell(9, 5).Range.Text = "86÷7="

$t.Cell(13, 1).Range.Text = "87÷5="
$t.Cell(13, 2).Range.Text = "90÷8="
$t.Cell(13, 3).Range.Text = "21÷5="
$t.Cell(13, 4).Range.Text = "62÷4="
$t.Cell(13, 5).Range.Text = "18÷5="

$t.Cell(17, 1).Range.Text = "54÷2="
$t.Cell(17, 2).Range.Text = "83÷9="
$t.Cell(17, 3).Range.Text = "71÷3="
$t.Cell(17, 4).Range.Text = "70÷3="
$t.Cell(17, 5).Range.Text = "89÷7="
